$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): columns E..M
#    Old layout: E=agri_land  F=precip  G=fert  H=pop  I=employ
#    New layout: E=agri_land_x F=agri_land_y G=precip H=fert_x I=fert_y
#                J=pop_x K=pop_y L=employ_x M=employ_y
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 5).Value = "4. Agriculture land area (% of land area)_x"
$ws.Cells.Item(1, 6).Value = "4. Agriculture land area (% of land area)_y"
$ws.Cells.Item(1, 7).Value = "5. Average precipitation (mm)"
$ws.Cells.Item(1, 8).Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_x"
$ws.Cells.Item(1, 9).Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_y"
$ws.Cells.Item(1, 10).Value = "13. Population_x"
$ws.Cells.Item(1, 11).Value = "13. Population_y"
$ws.Cells.Item(1, 12).Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_x"
$ws.Cells.Item(1, 13).Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_y"

# Give the brand-new header cells (J1:M1) the same formatting as the rest of
# the header row (bold font, border, centered) by copying format from A1.
$ws.Range("A1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Data rows. The sheet now carries five years of data (1997-2001) instead
#    of just 2000. Row 2 used to hold the (only) data row; we now rewrite it
#    and append four more rows below it, in ascending year order.
# ---------------------------------------------------------------------------

function Set-DataRow {
    param(
        [int]$Row,
        [int]$Year,
        [string]$Crop,
        [double]$Temp,
        [double]$Agri,
        [double]$Precip,
        [double]$Fert,
        [double]$Pop,
        [double]$Employ
    )

    $ws.Cells.Item($Row, 1).Value = "CPV"
    $ws.Cells.Item($Row, 2).Value = $Year

    # Column C ("0. Crop production index") keeps being stored as *text*
    # (matching the original workbook, where it is t="inlineStr" even
    # though the content looks numeric). Force text by temporarily
    # applying a text number format, then strip the format again so the
    # cell ends up without any explicit style, just like the source data.
    $cCell = $ws.Cells.Item($Row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $Crop
    $cCell.ClearFormats()

    $ws.Cells.Item($Row, 4).Value = $Temp
    $ws.Cells.Item($Row, 5).Value = $Agri
    $ws.Cells.Item($Row, 6).Value = $Agri
    $ws.Cells.Item($Row, 7).Value = $Precip
    $ws.Cells.Item($Row, 8).Value = $Fert
    $ws.Cells.Item($Row, 9).Value = $Fert
    $ws.Cells.Item($Row, 10).Value = $Pop
    $ws.Cells.Item($Row, 11).Value = $Pop
    $ws.Cells.Item($Row, 12).Value = $Employ
    $ws.Cells.Item($Row, 13).Value = $Employ
}

Set-DataRow 2 1997 "49.68" 23.09 17.86600496 125.37 2.318181818 430654 32.6455739613823
Set-DataRow 3 1998 "52.39" 23.14 17.86600496 145.98 2.590909091 440214 31.6164603587673
Set-DataRow 4 1999 "75.22" 23.15 17.86600496 148.15 4.25        449627 30.5989232839838
Set-DataRow 5 2000 "99.6"  23.27 17.86600496 118.98 2.931818182 458251 29.2854404606048
Set-DataRow 6 2001 "97.26" 23.42 18.36228288 144.98 5.282608696 465958 28.5888984158666

Write-Host "Edit complete. UsedRange:" $ws.UsedRange.Address()
